# HealthPassport - User Stories.docx
#
# The bullet "As a logged-in user, I want to be able to select multiple
# family members if a condition is applicable to more than one" is
# replaced by two bullets:
#   1) "As a logged-in user, I want to be able to add family history so
#       that I can get a full view of my medical risks stemming from my
#       family"   (same paragraph, split across 3 runs)
#   2) "As a logged-in user, I want to be edit family history so that I
#       can update my history whenever something happens to a family
#       member"   (new paragraph, same list style, split across 5 runs)

$d = $word.ActiveDocument

# Locate the paragraph to replace.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*select multiple family members*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

# Shared run properties (rFonts/color/sz/szCs) used by every run in both
# bullets -- identical to the run properties already on the paragraph
# being split.
$rPr = '<w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="2D3B45"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$xmlHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Rewrite the existing paragraph's text as 3 runs -----------------
$full = $target.Range
$bodyRange = $d.Range($full.Start, $full.End - 1)

$para1Runs = `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">As a logged-in user, I want to be able to </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>add family history</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> so that I can get a full view of my medical risks stemming from my family</w:t></w:r>'

$para1Xml = $xmlHeader + '<w:p>' + $para1Runs + '</w:p>' + $xmlFooter
$bodyRange.InsertXML($para1Xml)

# --- Insert a brand-new paragraph right after it, same list style ----
$pPr = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="180" w:after="180" w:line="240" w:lineRule="auto"/>' + $rPr + '</w:pPr>'

$para2Runs = `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">As a logged-in user, I want to be </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>edit</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> family history</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> so that I can</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> update my history whenever something happens to a family member</w:t></w:r>'

# Create an empty paragraph right after the (now rewritten) first bullet,
# then fill its whole range (including the paragraph mark) with the new
# paragraph + a fresh trailing paragraph mark so nothing gets merged
# into the paragraph that follows.
$null = $target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newRange = $newPara.Range

$para2Xml = $xmlHeader + '<w:p>' + $pPr + $para2Runs + '</w:p><w:p/>' + $xmlFooter
$insertRange = $d.Range($newRange.Start, $newRange.End)
$insertRange.InsertXML($para2Xml)
